$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1060579828356883
$ws.Range("E2").Value = 0.2098744004587479
$ws.Range("F2").Value = 0.1841535432926422

$ws.Range("D3").Value = 0.5257115046245979
$ws.Range("E3").Value = 0.5521715976014016
$ws.Range("F3").Value = 0.863903320853156

$ws.Range("D4").Value = 0.2269096954953032
$ws.Range("E4").Value = 0.4491329609791932
$ws.Range("F4").Value = 1

$ws.Range("D5").Value = 0.2507116742946722
$ws.Range("E5").Value = 0.3445457345191699
$ws.Range("F5").Value = 0.5941651397226838

$ws.Range("D6").Value = 0.2497033435456695
$ws.Range("E6").Value = 0.6211223257768826
$ws.Range("F6").Value = 0.6302434217587324

$ws.Range("D7").Value = 0.09100679477754352
$ws.Range("E7").Value = 0.1322881038858737
$ws.Range("F7").Value = 0.9059275295124692

$ws.Range("D8").Value = 0.05379139958613055
$ws.Range("E8").Value = 0.0853604963672745
$ws.Range("F8").Value = 1

$ws.Range("D9").Value = 0.1007474296530028
$ws.Range("E9").Value = 0.1323748351993976
$ws.Range("F9").Value = 0.1176378875462932

$ws.Range("D10").Value = 0.201519929683999
$ws.Range("E10").Value = 0.2318872633403544
$ws.Range("F10").Value = 0.617829656649348

$ws.Range("D11").Value = 0.2350631375161907
$ws.Range("E11").Value = 0.2610587746687676
$ws.Range("F11").Value = 0.4329073271576331

$ws.Range("D15").Value = 0.205141630736422
$ws.Range("E15").Value = 0.3626649014229366
$ws.Range("F15").Value = 0.9413985616468207
